$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2").Value = 'buluş, millet, sandık, referandum, seçim'
$ws.Range("B2").Value = 68
$ws.Range("C2").Value = 'millet, oy, buluş, sandık, seçim'
$ws.Range("A3").Value = 'atatürk, önder, kemal, kutlu, gazi'
$ws.Range("B3").Value = 16
$ws.Range("C3").Value = 'bayram, atatürk, türk, kutlu, mustafa'
$ws.Range("A4").Value = 'tv, program, konuk, yayın, fox'
$ws.Range("B4").Value = 14
$ws.Range("C4").Value = 'yayın, program, canlı, konuk, sun'
$ws.Range("A5").Value = 'sağduyu, pis, yeter, çık, allah'
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 'iyi, sanatçı, allah, iş, çık'
$ws.Range("A6").Value = 'basın, açıkla, kktc, dön, medya'
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 'basın, açıkla, medya, uygula, cemiyet'
$ws.Range("A7").Value = 'milyar, yatırım, para, dolar, lira'
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = 'milyar, yatırım, lira, dolar, kamu'
$ws.Range("A8").Value = 'türkiye, türk, dostluk, israil, hemen'
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = 'türk, türkiye, milliyetçi, yüzyıl, millet'
$ws.Range("A9").Value = 'bura, gel, haydi, buluş, bitir'
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 'gel, buluş, saat, bugün, bekle'
$ws.Range("A10").Value = 'lig, süper, rahmetli, camia, ateş'
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 'şehit, rahmet, atatürk, an, dönüm'
$ws.Range("A11").Value = 'deprem, depremzede, konut, bölge, hele'
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = 'konut, deprem, depremzede, temel, hastane'
$ws.Range("A12").Value = 'genç, demirel, ak, salon, partili'
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = 'genç, internet, buluş, telefon, medya'
$ws.Range("A13").Value = 'acı, kayıp, çerkes, sürgün, din'
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = 'acı, çerkes, sürgün, kardeş, kayıp'
